$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.610.96"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.923.07"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.32"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2890"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06836"
$ws.Range("E9").Value = "  +3.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "105.12"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.37"
$ws.Range("E11").Value = "  -4.08%  "
$ws.Range("D12").Value = "1.915.97"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07699"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("E14").Value = "  +3.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6678"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "292.19"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").Value = "30.608.08"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.596"
$ws.Range("E18").Value = "  +6.71%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007624"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.96"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "2.176.78"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.430"
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.428"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.82"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.08"
$ws.Range("E27").Value = "  +7.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.119"
$ws.Range("E28").Value = "  +4.66%  "
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.181"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("E32").Value = "  +3.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05046"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7382"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.145"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02078"
$ws.Range("E36").Value = "  +6.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.738"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.692"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.059"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "111.24"
$ws.Range("E40").Value = "  +3.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8755"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4392"
$ws.Range("E42").Value = "  +6.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.888"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.05"
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.270"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.373"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.36"
$ws.Range("E48").Value = "  +15.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1244"
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.2520"
$ws.Range("E50").Value = "  +13.33%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.03"
$ws.Range("E51").Value = "  +1.16%  "
